$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Bottling day" row: the recorded date (A7) moves from 2020-11-03 to
# 2020-11-04 (serial 44138 -> 44139). B7 (=A7-$A$6, the "Day #" column)
# depends on A7 and recalculates automatically (10 -> 11).
$ws.Range("A7").Value = 44139
